$d = $word.ActiveDocument

# The new bullet point is added right after the existing last paragraph
# ("We could use bar graphs ... concentrated.") and before the section
# properties. InsertParagraphAfter() on that paragraph's range creates a
# brand-new paragraph that inherits the same paragraph formatting
# (ListParagraph style + the numbered/bulleted list numPr), which matches
# the <w:pPr> block shown in the diff for the new paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

# Grab the freshly created (currently empty) paragraph and fill it in with
# the new sentence.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = "The median would be better to summarize the data because of the high variability of the data would make the mean significantly higher than where a large portion of the data is concentrated."
